$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename header cells: "..._old" -> "..._FV2304" (the earlier format version),
#    "..._new" -> "..._FV2310" (the later format version).
for ($c = 1; $c -le 21; $c++) {
    $cell = $ws.Cells.Item(1, $c)
    $val = $cell.Value()
    if ($val -ne $null) {
        if ($val.EndsWith("_old")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2304"
        } elseif ($val.EndsWith("_new")) {
            $cell.Value = $val.Substring(0, $val.Length - 4) + "_FV2310"
        }
    }
}

# 2. Turn the used range into a real Excel Table ("Table1") with an AutoFilter,
#    picking up the freshly renamed header row as its column names.
$range = $ws.Range("A1:U56")
$tbl = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $range, $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

# 3. Freeze the header row (split below row 1, top-left of scrolling area = A2).
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
